$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2735
$ws.Range("J2").Value = 4900
$ws.Range("L2").Value = 4900
$ws.Range("N2").Value = -5126
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216
$ws.Range("H19").Value = 867.95
$ws.Range("I19").Value = 932
$ws.Range("J19").Value = 771.875
$ws.Range("K19").Value = 932
$ws.Range("L19").Value = 771.875
$ws.Range("M19").Value = -757
$ws.Range("N19").Value = -1121.875
$ws.Range("H40").Value = 3305.5
$ws.Range("I40").Value = 3407.3333
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3407.3333
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3232.3333
$ws.Range("N40").Value = -3350
$ws.Range("H74").Value = 168966.33
$ws.Range("I74").Value = 3449.5
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 3449.5
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -2513.5
$ws.Range("N74").Value = -501872
$ws.Range("H77").Value = 168966.33
$ws.Range("I77").Value = 3449.5
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 17247.5
$ws.Range("L77").Value = 2500000
$ws.Range("M77").Value = -12567.5
$ws.Range("N77").Value = -2509360
$ws.Range("H107").Value = 630
$ws.Range("I107").Value = 561.875
$ws.Range("K107").Value = 561.875
$ws.Range("M107").Value = 1358.125
$ws.Range("H111").Value = 947
$ws.Range("I111").Value = 899.5
$ws.Range("K111").Value = 2698.5
$ws.Range("M111").Value = 368.5
$ws.Range("H118").Value = 1299
$ws.Range("I118").Value = 1299
$ws.Range("K118").Value = 3897
$ws.Range("M118").Value = -2240
$ws.Range("H127").Value = 2099.1667
$ws.Range("I127").Value = 2065
$ws.Range("J127").Value = 2133.3333
$ws.Range("K127").Value = 6195
$ws.Range("L127").Value = 6399.999899999999
$ws.Range("M127").Value = -1235
$ws.Range("N127").Value = -16319.9999
$ws.Range("H129").Value = 1768.4
$ws.Range("I129").Value = 615.1667
$ws.Range("J129").Value = 3498.25
$ws.Range("K129").Value = 1845.5001
$ws.Range("L129").Value = 10494.75
$ws.Range("M129").Value = 3154.4999
$ws.Range("N129").Value = -20494.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1875
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1123
$ws.Range("H61").Value = 1900
$ws.Range("I61").Value = 1900
$ws.Range("K61").Value = 1900
$ws.Range("M61").Value = -1688
$ws.Range("H110").Value = 2437
$ws.Range("I110").Value = 1412.4
$ws.Range("K110").Value = 1412.4
$ws.Range("M110").Value = 632.5999999999999
$ws.Range("H122").Value = 2838.5715
$ws.Range("I122").Value = 1493.5
$ws.Range("K122").Value = 4480.5
$ws.Range("M122").Value = -2030.5
$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -3150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3403.1667
$ws.Range("I105").Value = 3173.8
$ws.Range("J105").Value = 4550
$ws.Range("K105").Value = 3173.8
$ws.Range("L105").Value = 4550
$ws.Range("M105").Value = -1426.8
$ws.Range("N105").Value = -8044

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1797.1111
$ws.Range("I31").Value = 1797.1111
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1797.1111
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1502.1111
$ws.Range("N31").Value = $null
$ws.Range("H34").Value = 1797.1111
$ws.Range("I34").Value = 1797.1111
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1797.1111
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1595.1111
$ws.Range("N34").Value = $null
$ws.Range("H35").Value = 717460.5600000001
$ws.Range("I35").Value = 717460.5600000001
$ws.Range("K35").Value = 717460.5600000001
$ws.Range("M35").Value = -717166.5600000001
$ws.Range("H107").Value = 469.4
$ws.Range("I107").Value = 426.125
$ws.Range("K107").Value = 426.125
$ws.Range("M107").Value = 1493.875
$ws.Range("H132").Value = 2851.1875
$ws.Range("I132").Value = 1924.3334
$ws.Range("K132").Value = 5773.0002
$ws.Range("M132").Value = -3243.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163
$ws.Range("I2").Value = 231.25
$ws.Range("J2").Value = 26.5
$ws.Range("K2").Value = 1387.5
$ws.Range("L2").Value = 159
$ws.Range("M2").Value = -1274.5
$ws.Range("N2").Value = -385
$ws.Range("H5").Value = 829.2
$ws.Range("I5").Value = 788.25
$ws.Range("J5").Value = 993
$ws.Range("K5").Value = 2364.75
$ws.Range("L5").Value = 2979
$ws.Range("M5").Value = -2252.75
$ws.Range("N5").Value = -3203
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H98").Value = 151.5
$ws.Range("I98").Value = 110
$ws.Range("J98").Value = 165.33333
$ws.Range("K98").Value = 330
$ws.Range("L98").Value = 495.99999
$ws.Range("M98").Value = 1168
$ws.Range("N98").Value = -3491.99999
$ws.Range("H135").Value = 829.2
$ws.Range("I135").Value = 788.25
$ws.Range("J135").Value = 993
$ws.Range("K135").Value = 7094.25
$ws.Range("L135").Value = 8937
$ws.Range("M135").Value = -4559.25
$ws.Range("N135").Value = -14007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3750
$ws.Range("I80").Value = 3500
$ws.Range("K80").Value = 3500
$ws.Range("M80").Value = -2502
$ws.Range("H83").Value = 3750
$ws.Range("I83").Value = 3500
$ws.Range("K83").Value = 17500
$ws.Range("M83").Value = -12508
$ws.Range("H87").Value = 19000
$ws.Range("I87").Value = 19000
$ws.Range("K87").Value = 19000
$ws.Range("M87").Value = -17752
$ws.Range("H90").Value = 19000
$ws.Range("I90").Value = 19000
$ws.Range("K90").Value = 57000
$ws.Range("M90").Value = -50760
$ws.Range("H107").Value = 175.2
$ws.Range("I107").Value = 182
$ws.Range("K107").Value = 182
$ws.Range("M107").Value = 1738
$ws.Range("H122").Value = 1562
$ws.Range("I122").Value = 1556.25
$ws.Range("K122").Value = 4668.75
$ws.Range("M122").Value = -2218.75
$ws.Range("H126").Value = 11488
$ws.Range("I126").Value = 11488
$ws.Range("K126").Value = 34464
$ws.Range("M126").Value = -31994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = $null
$ws.Range("H40").Value = 7378.625
$ws.Range("I40").Value = 7378.625
$ws.Range("K40").Value = 7378.625
$ws.Range("M40").Value = -7242.625
$ws.Range("H120").Value = 19999
$ws.Range("J120").Value = 19999
$ws.Range("L120").Value = 19999
$ws.Range("N120").Value = -29675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10000
$ws.Range("H107").Value = 2897.4666
$ws.Range("I107").Value = 3386.4443
$ws.Range("K107").Value = 10159.3329
$ws.Range("M107").Value = -8239.332900000001
$ws.Range("H113").Value = 288.6
$ws.Range("I113").Value = 245.66667
$ws.Range("K113").Value = 737.00001
$ws.Range("M113").Value = 1432.99999
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
